$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "content"
$ws.Range("C1").Value = "tag"
$ws.Range("D1").Value = "label"
$ws.Range("E1").Value = "link"
$ws.Range("F1").Value = "output"
$ws.Range("G1").Value = "prompt version"

# --- Data row (row 2) ---
$ws.Range("A2").Value = "莎普爱思滴眼睛"

$content = "莎普爱思滴眼睛`t用词模糊,消息虚假,夸大失实`t白内障，看不清~莎普爱思滴眼睛~`n白内障，看不清~莎普爱思滴滴滴！`n模糊滴！重影滴！黑影滴！`n白内障，看不清~莎普爱思滴眼睛~`n白内障，看不清~莎普爱思滴滴滴！`n模糊滴！重影滴！黑影滴！有点痛！坚持滴！！！。 "
$ws.Range("B2").Value = $content

$ws.Range("C2").Value = "用词模糊,消息虚假,夸大失实"

# D2 and E2 (label / link) are intentionally left blank

$output = "我帮你看了一下哈～这个广告使用了模糊的用词和夸大的手法，来宣传莎普爱思滴眼液对白内障的治疗效果。虽然滴眼液可能对某些人有帮助，但这种夸大和不清晰的表述可能存在误导性。建议你谨慎对待这样的广告宣传，并且在购买前咨询专业的医生或药剂师，获取更准确的产品信息和使用建议。避免购买🚫"
$ws.Range("F2").Value = $output

$ws.Range("G2").Value = "v2"
